# edit.ps1 - Apply the "Wed May 3 18:08:30 UTC 2023" cryptos-list refresh.
# For each changed row: updates the Price (col D) and/or Volume(1h) (col E) text,
# and for rows 45/46 also updates Coin (B) and Link (C), since the two rows swap
# which coin (EnergySwap / Decentraland) occupies which position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet stores Price/Volume as plain text (not numbers), e.g. "87.00" or
# "0.00001000" - values whose trailing/leading zeros must survive verbatim and
# which use "." as a thousands separator in a few rows (e.g. "28.737.19").
# Plain `.Value = "..."` assignment lets Excel auto-detect such look-alike
# numeric strings and coerce them to real numbers (dropping the formatting we
# need, e.g. "87.00" -> 87). Pre-formatting the target cells as Text keeps the
# literal string, matching the source data.
$textCells = @("D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D19", "D20", "D23", "D24", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.737.19"
$ws.Range("E2").Value = "  +0.12%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.890.97"
$ws.Range("E3").Value = "  +1.21%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.04%  "

# Row 5 - BNB
$ws.Range("D5").Value = "326.16"
$ws.Range("E5").Value = "  -0.25%  "

# Row 6 - USDC
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.16%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.4568"
$ws.Range("E7").Value = "  -1.41%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.3856"
$ws.Range("E8").Value = "  -1.45%  "

# Row 9 - OKB
$ws.Range("D9").Value = "46.63"
$ws.Range("E9").Value = "  +0.45%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.07862"
$ws.Range("E10").Value = "  -0.58%  "

# Row 11 - Polygon
$ws.Range("D11").Value = "0.9996"
$ws.Range("E11").Value = "  +3.12%  "

# Row 12 - Solana
$ws.Range("D12").Value = "21.68"
$ws.Range("E12").Value = "  -2.73%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.908.87"
$ws.Range("E13").Value = "  +3.50%  "

# Row 14 - Chainlink
$ws.Range("D14").Value = "6.990"
$ws.Range("E14").Value = "  +0.88%  "

# Row 15 - Polkadot
$ws.Range("D15").Value = "5.688"
$ws.Range("E15").Value = "  -0.77%  "

# Row 16 - TRON
$ws.Range("D16").Value = "0.06950"
$ws.Range("E16").Value = "  -0.15%  "

# Row 17 - Litecoin
$ws.Range("D17").Value = "87.00"
$ws.Range("E17").Value = "  -1.40%  "

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.00001000"
$ws.Range("E19").Value = "  -0.57%  "

# Row 20 - Avalanche
$ws.Range("D20").Value = "16.90"
$ws.Range("E20").Value = "  -0.27%  "

# Row 22 - WrappedBTC
$ws.Range("D22").Value = "28.767.52"
$ws.Range("E22").Value = "  +0.24%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "5.280"
$ws.Range("E23").Value = "  -0.77%  "

# Row 24 - Cosmos
$ws.Range("D24").Value = "10.92"
$ws.Range("E24").Value = "  -1.48%  "

# Row 25 - WrappedliquidstakedEther2.0
$ws.Range("D25").Value = "2.141.27"
$ws.Range("E25").Value = "  +1.29%  "

# Row 26 - Toncoin
$ws.Range("D26").Value = "2.069"
$ws.Range("E26").Value = "  -2.49%  "

# Row 27 - Monero
$ws.Range("D27").Value = "154.55"
$ws.Range("E27").Value = "  +0.71%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  -0.74%  "

# Row 29 - InternetComputer(DFINITY)
$ws.Range("D29").Value = "5.781"
$ws.Range("E29").Value = "  +1.26%  "

# Row 30 - BitcoinCash
$ws.Range("D30").Value = "118.27"
$ws.Range("E30").Value = "  -1.04%  "

# Row 31 - LidoDAOToken
$ws.Range("D31").Value = "1.903"
$ws.Range("E31").Value = "  -4.82%  "

# Row 32 - Stellar
$ws.Range("D32").Value = "0.09294"
$ws.Range("E32").Value = "  -0.72%  "

# Row 33 - ImmutableX
$ws.Range("D33").Value = "0.9141"
$ws.Range("E33").Value = "  -1.82%  "

# Row 34 - Filecoin
$ws.Range("E34").Value = "  -0.47%  "

# Row 35 - ARBITRUM
$ws.Range("D35").Value = "1.327"
$ws.Range("E35").Value = "  -1.33%  "

# Row 36 - HuobiToken
$ws.Range("D36").Value = "3.259"
$ws.Range("E36").Value = "  -2.93%  "

# Row 37 - Hedera
$ws.Range("D37").Value = "0.05708"
$ws.Range("E37").Value = "  -2.14%  "

# Row 38 - TrustWalletToken
$ws.Range("D38").Value = "1.155"
$ws.Range("E38").Value = "  +0.56%  "

# Row 39 - VeChain
$ws.Range("D39").Value = "0.02055"
$ws.Range("E39").Value = "  -3.48%  "

# Row 40 - FraxShare
$ws.Range("D40").Value = "7.673"
$ws.Range("E40").Value = "  -2.95%  "

# Row 41 - TheSandbox
$ws.Range("D41").Value = "0.5601"
$ws.Range("E41").Value = "  -0.96%  "

# Row 42 - Algorand
$ws.Range("E42").Value = "  -0.33%  "

# Row 43 - Aptos
$ws.Range("D43").Value = "9.696"
$ws.Range("E43").Value = "  -2.63%  "

# Row 44 - Cronos
$ws.Range("D44").Value = "0.07155"
$ws.Range("E44").Value = "  -1.07%  "

# Row 45 - EnergySwap
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "0.5290"
$ws.Range("E45").Value = "  -0.41%  "

# Row 46 - Decentraland
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "11.58"
$ws.Range("E46").Value = "  -1.71%  "

# Row 47 - RenderToken
$ws.Range("D47").Value = "2.159"
$ws.Range("E47").Value = "  +1.31%  "

# Row 48 - WEMIXToken
$ws.Range("D48").Value = "1.119"
$ws.Range("E48").Value = "  -1.49%  "

# Row 49 - NEARProtocol
$ws.Range("D49").Value = "1.819"
$ws.Range("E49").Value = "  -1.54%  "

# Row 50 - Quant
$ws.Range("D50").Value = "112.04"
$ws.Range("E50").Value = "  -1.26%  "

# Row 51 - MXToken
$ws.Range("D51").Value = "2.454"
$ws.Range("E51").Value = "  +4.47%  "
